$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-02-29 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-01 Friday", 2) | Out-Null

# Update the division-problem answers inside the table, cell by cell
# (addressed positionally since some cell values repeat verbatim)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "94÷6=15, 4"
$t.Cell(1, 2).Range.Text = "50÷9=5, 5"
$t.Cell(1, 3).Range.Text = "75÷8=9, 3"
$t.Cell(1, 4).Range.Text = "81÷9=9, 0"
$t.Cell(1, 5).Range.Text = "55÷6=9, 1"
$t.Cell(5, 1).Range.Text = "56÷9=6, 2"
$t.Cell(5, 2).Range.Text = "43÷2=21, 1"
$t.Cell(5, 3).Range.Text = "49÷9=5, 4"
$t.Cell(5, 4).Range.Text = "56÷9=6, 2"
$t.Cell(5, 5).Range.Text = "66÷4=16, 2"
$t.Cell(9, 1).Range.Text = "33÷8=4, 1"
$t.Cell(9, 2).Range.Text = "53÷3=17, 2"
$t.Cell(9, 3).Range.Text = "78÷3=26, 0"
$t.Cell(9, 4).Range.Text = "99÷3=33, 0"
$t.Cell(9, 5).Range.Text = "25÷2=12, 1"
$t.Cell(13, 1).Range.Text = "44÷2=22, 0"
$t.Cell(13, 2).Range.Text = "71÷5=14, 1"
$t.Cell(13, 3).Range.Text = "77÷2=38, 1"
$t.Cell(13, 4).Range.Text = "41÷9=4, 5"
$t.Cell(13, 5).Range.Text = "12÷7=1, 5"
$t.Cell(17, 1).Range.Text = "60÷8=7, 4"
$t.Cell(17, 2).Range.Text = "98÷7=14, 0"
$t.Cell(17, 3).Range.Text = "13÷3=4, 1"
$t.Cell(17, 4).Range.Text = "43÷6=7, 1"
$t.Cell(17, 5).Range.Text = "50÷2=25, 0"
